{"js": "// Rewrite the \"KEY ACHIEVEMENTS AND IMPACT\" bullet list (under the \"Impact\"\n// Heading3) from 6 job-duty-style bullets into 4 impact-focused accomplishment\n// statements, per the commit:\n//   \"Fix Key Achievements to use proper accomplishment statements\"\n//\n// Before (6 bullets):\n//   1. Built redistricting platform used by thousands of analysts nationwide\n//      with real-time collaborative editing and Census integration, serving\n//      12,847 analysts across 89 organizations\n//   2. Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for\n//      large-scale geospatial datasets\n//   3. Trigonometric algorithm for boundary estimation reduced mapping costs\n//      by 73.5%, saving campaigns and organizations $4.7M and enabling\n//      smaller nonprofits to conduct analysis\n//   4. Discovered systematic race coding errors affecting all Black and\n//      Asian-American voters, developed geospatial machine learning\n//      algorithms improving classification accuracy from 23% to 64%\n//   5. Achieved 87% prediction accuracy for voter turnout vs. industry\n//      standard of 71%, reducing polling error margins from \u00b14.2% to \u00b12.1%\n//   6. Built cloud-based data warehouse solutions on AWS processing billions\n//      of records with 99.94% accuracy\n//\n// After (4 bullets):\n//   1. Algorithmic innovation: Pioneered trigonometric boundary estimation\n//      reducing mapping costs 73.5%\n//   2. $4.7M savings enabled nonprofit access\n//   3. Platform impact: Built redistricting system serving 12,847 analysts\n//      across 89 organizations\n//   4. Real-time collaboration at national scale\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the \"KEY ACHIEVEMENTS AND IMPACT\" heading, then its \"Impact\"\n// sub-heading directly after it, then operate on the bullet paragraphs that\n// immediately follow -- this avoids accidentally touching the near-duplicate\n// bullet text that also appears earlier under \"Partner - Siege Analytics\".\nconst items = paragraphs.items;\nlet achIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.trim() === \"KEY ACHIEVEMENTS AND IMPACT\") {\n    achIdx = i;\n    break;\n  }\n}\nif (achIdx === -1) {\n  throw new Error('Could not find \"KEY ACHIEVEMENTS AND IMPACT\" heading');\n}\n\nlet impactIdx = -1;\nfor (let i = achIdx + 1; i < items.length; i++) {\n  if (items[i].text.trim() === \"Impact\") {\n    impactIdx = i;\n    break;\n  }\n}\nif (impactIdx === -1) {\n  throw new Error('Could not find \"Impact\" sub-heading after KEY ACHIEVEMENTS AND IMPACT');\n}\n\n// The six bullet paragraphs follow immediately after the \"Impact\" heading.\nconst bulletStart = impactIdx + 1;\nconst oldBullets = [\n  \"\u2022 Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations\",\n  \"\u2022 Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets\",\n  \"\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis\",\n  \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%\",\n  \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \u00b14.2% to \u00b12.1%\",\n  \"\u2022 Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy\",\n];\n\nfor (let i = 0; i < oldBullets.length; i++) {\n  const actual = items[bulletStart + i].text.trim();\n  if (actual !== oldBullets[i]) {\n    throw new Error(\n      \"Unexpected bullet text at index \" + (bulletStart + i) + \": \" + actual\n    );\n  }\n}\n\nconst newBullets = [\n  \"\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\",\n  \"\u2022 $4.7M savings enabled nonprofit access\",\n  \"\u2022 Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations\",\n  \"\u2022 Real-time collaboration at national scale\",\n];\n\n// Replace the text of the first four bullet paragraphs with the new content.\nfor (let i = 0; i < newBullets.length; i++) {\n  const para = items[bulletStart + i];\n  para.getRange().insertText(newBullets[i], Word.InsertLocation.replace);\n}\n\n// Delete the two now-unneeded trailing bullet paragraphs (old bullets 5 & 6,\n// i.e. indices bulletStart+4 and bulletStart+5). Delete from the end first so\n// indices stay valid.\nitems[bulletStart + 5].delete();\nitems[bulletStart + 4].delete();\n\nawait context.sync();\n", "ps1": "# Rewrite the \"KEY ACHIEVEMENTS AND IMPACT\" bullet list (under the \"Impact\"\n# Heading3) from 6 job-duty-style bullets into 4 impact-focused accomplishment\n# statements, per the commit:\n#   \"Fix Key Achievements to use proper accomplishment statements\"\n#\n# Before (6 bullets):\n#   1. Built redistricting platform used by thousands of analysts nationwide\n#      with real-time collaborative editing and Census integration, serving\n#      12,847 analysts across 89 organizations\n#   2. Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for\n#      large-scale geospatial datasets\n#   3. Trigonometric algorithm for boundary estimation reduced mapping costs\n#      by 73.5%, saving campaigns and organizations $4.7M and enabling\n#      smaller nonprofits to conduct analysis\n#   4. Discovered systematic race coding errors affecting all Black and\n#      Asian-American voters, developed geospatial machine learning\n#      algorithms improving classification accuracy from 23% to 64%\n#   5. Achieved 87% prediction accuracy for voter turnout vs. industry\n#      standard of 71%, reducing polling error margins from \u00b14.2% to \u00b12.1%\n#   6. Built cloud-based data warehouse solutions on AWS processing billions\n#      of records with 99.94% accuracy\n#\n# After (4 bullets):\n#   1. Algorithmic innovation: Pioneered trigonometric boundary estimation\n#      reducing mapping costs 73.5%\n#   2. $4.7M savings enabled nonprofit access\n#   3. Platform impact: Built redistricting system serving 12,847 analysts\n#      across 89 organizations\n#   4. Real-time collaboration at national scale\n\n$d = $word.ActiveDocument\n\nfunction Get-ParaText($idx) {\n    $t = $d.Paragraphs.Item($idx).Range.Text\n    return $t.TrimEnd([char]13, [char]7, [char]10)\n}\n\n# Locate the \"KEY ACHIEVEMENTS AND IMPACT\" heading, then its \"Impact\"\n# sub-heading directly after it, then operate on the bullet paragraphs that\n# immediately follow -- this avoids accidentally touching the near-duplicate\n# bullet text that also appears earlier under \"Partner - Siege Analytics\".\n$achIdx = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ((Get-ParaText $i) -eq \"KEY ACHIEVEMENTS AND IMPACT\") {\n        $achIdx = $i\n        break\n    }\n}\nif ($achIdx -eq -1) { throw \"Could not find 'KEY ACHIEVEMENTS AND IMPACT' heading\" }\n\n$impactIdx = -1\nfor ($i = $achIdx + 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ((Get-ParaText $i) -eq \"Impact\") {\n        $impactIdx = $i\n        break\n    }\n}\nif ($impactIdx -eq -1) { throw \"Could not find 'Impact' sub-heading after KEY ACHIEVEMENTS AND IMPACT\" }\n\n# The six bullet paragraphs follow immediately after the \"Impact\" heading.\n$bulletStart = $impactIdx + 1\n\n$oldBullets = @(\n    \"\u2022 Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations\",\n    \"\u2022 Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets\",\n    \"\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis\",\n    \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%\",\n    \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \u00b14.2% to \u00b12.1%\",\n    \"\u2022 Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy\"\n)\n\nfor ($k = 0; $k -lt $oldBullets.Length; $k++) {\n    $actual = Get-ParaText ($bulletStart + $k)\n    if ($actual -ne $oldBullets[$k]) {\n        throw \"Unexpected bullet text at paragraph $($bulletStart + $k): $actual\"\n    }\n}\n\n$newBullets = @(\n    \"\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\",\n    \"\u2022 $4.7M savings enabled nonprofit access\",\n    \"\u2022 Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations\",\n    \"\u2022 Real-time collaboration at national scale\"\n)\n\n# Replace the text of the first four bullet paragraphs with the new content.\nfor ($k = 0; $k -lt $newBullets.Length; $k++) {\n    $d.Paragraphs.Item($bulletStart + $k).Range.Text = $newBullets[$k]\n}\n\n# Delete the two now-unneeded trailing bullet paragraphs (old bullets 5 & 6,\n# i.e. paragraphs bulletStart+4 and bulletStart+5). Delete from the end first\n# so indices stay valid.\n$d.Paragraphs.Item($bulletStart + 5).Range.Delete()\n$d.Paragraphs.Item($bulletStart + 4).Range.Delete()\n"}
